# Apply cryptos list update (prices & 1h volume % changes) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.597.77'
$ws.Range('E2').Value = '  -2.95%  '
$ws.Range('D3').Value = '1.980.44'
$ws.Range('E3').Value = '  -3.86%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'246.12"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.00%  '
$ws.Range('E6').Value = '  -5.01%  '
$ws.Range('D7').Value = "'58.59"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +7.34%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = "'58.73"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.20%  '
$ws.Range('E10').Value = '  -1.21%  '
$ws.Range('D11').Value = "'0.0735"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -2.20%  '
$ws.Range('E12').Value = '  -2.81%  '
$ws.Range('D13').Value = "'0.938"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.69%  '
$ws.Range('D14').Value = "'14.51"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.59%  '
$ws.Range('D15').Value = '2.271.14'
$ws.Range('E15').Value = '  -3.84%  '
$ws.Range('E16').Value = '  -3.09%  '
$ws.Range('D17').Value = '1.991.03'
$ws.Range('E17').Value = '  -3.47%  '
$ws.Range('D18').Value = "'18.05"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +5.87%  '
$ws.Range('D19').Value = '35.514.14'
$ws.Range('E19').Value = '  -2.99%  '
$ws.Range('D20').Value = "'71.39"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.99%  '
$ws.Range('D21').Value = '0.0₃0846'
$ws.Range('E21').Value = '  -2.19%  '
$ws.Range('D22').Value = "'5.21"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.09%  '
$ws.Range('D23').Value = "'232.65"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.48%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').Value = "'2.63"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +23.22%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').Value = "'1.00"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('D26').Value = "'2.28"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -3.83%  '
$ws.Range('D27').Value = "'164.83"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.23%  '
$ws.Range('D28').Value = "'9.12"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.55%  '
$ws.Range('E29').Value = '  -5.01%  '
$ws.Range('E30').Value = '  -2.72%  '
$ws.Range('D31').Value = "'4.85"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -4.76%  '
$ws.Range('E32').Value = '  -6.64%  '
$ws.Range('D33').Value = "'0.0970"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +17.54%  '
$ws.Range('E34').Value = '  -0.62%  '
$ws.Range('D35').Value = "'2.42"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +9.20%  '
$ws.Range('E36').Value = '  -3.57%  '
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('E38').Value = '  -3.84%  '
$ws.Range('D39').Value = "'5.37"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +9.49%  '
$ws.Range('E40').Value = '  -2.45%  '
$ws.Range('E41').Value = '  -0.62%  '
$ws.Range('E42').Value = '  -1.92%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = "'7.80"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +2.81%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = "'93.70"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.80%  '
$ws.Range('E45').Value = '  -1.71%  '
$ws.Range('D46').Value = "'16.15"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.92%  '
$ws.Range('D47').Value = "'0.0895"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.91%  '
$ws.Range('D48').Value = '1.374.04'
$ws.Range('E48').Value = '  -2.38%  '
$ws.Range('E49').Value = '  -0.64%  '
$ws.Range('D50').Value = "'47.18"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.51%  '
$ws.Range('E51').Value = '  -0.30%  '

Write-Host "Applied crypto list update"
